# perDW.xlsx edit: add a "Sum" column (S) that totals the per-compound
# columns (B:R) on each data sheet, and switch the active sheet/selection
# back to "BA".

$wb = $excel.ActiveWorkbook

$sheetNames = @("BA", "N")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # header for the new column
    $ws.Cells.Item(1, 19).Value = "Sum"

    $lastRow = $ws.UsedRange.Rows.Count

    for ($row = 2; $row -le $lastRow; $row++) {
        $total = 0.0
        for ($col = 2; $col -le 18; $col++) {
            $total = $total + $ws.Cells.Item($row, $col).Value()
        }
        $ws.Cells.Item($row, 19).Value = $total
    }
}

# Cosmetic print-setup tweak that came along with this edit on "N".
try {
    $wb.Worksheets.Item("N").PageSetup.FirstPageNumber = 0
} catch {
}

# Restore "BA" as the active sheet/tab, with the selection sitting on the
# new Sum header (R1 -> now logically S1 region), and leave "N" selected
# at U5 (its last-used selection) but no longer the active tab.
$wsN = $wb.Worksheets.Item("N")
$wsN.Activate()
$wsN.Range("U5").Select()

$wsBA = $wb.Worksheets.Item("BA")
$wsBA.Activate()
$wsBA.Range("R1").Select()
